# Apply the "Update resolved points" revision to the SRS review workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Introduction " ---
$ws1 = $wb.Worksheets.Item(1)

# Ref Version bumped from 1 to 0.6
$ws1.Range("D7").Value = 0.6

# Last update date bumped to 23/2/2020
$ws1.Range("D9").Value = "23/2/2020"

# New history row documenting the resolved-points update
$ws1.Range("B14").Value = 0.2
$ws1.Range("C14").Value = "T.Sharaby"
$ws1.Range("E14").Value = "23/2/2020"
$ws1.Range("G14").Value = "Update resolved points"

# --- Sheet 2: "Cross review points " ---
$ws2 = $wb.Worksheets.Item(2)

# Mark the open review points as resolved
$ws2.Range("H2:H7").Value = "Resolved"

# Clear the status for the still-blank row
$ws2.Range("H8").Value = ""
